# Auto update Excel log
# Appends new sensor-log rows to the ALERTS sheet and the mmWave sheet,
# matching the data produced by the logging process on 2026-02-01.

$wb = $excel.ActiveWorkbook

function Add-LogRow($Sheet, $Row, $DateVal, $TimeVal, $HourVal, $LocationVal, $ValueVal, $StatusVal) {
    # Column A holds date-looking text (e.g. "2026-02-01"). Excel would
    # normally auto-convert such text into a real date serial number, so
    # force the cell to Text format first, then reset the style back to
    # Normal/General once the literal text has been stored - this keeps
    # the stored value as plain text without leaving a custom number
    # format behind.
    $cellA = $Sheet.Cells.Item($Row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $DateVal
    $cellA.Style = "Normal"

    $Sheet.Cells.Item($Row, 2).Value = $TimeVal
    $Sheet.Cells.Item($Row, 3).Value = $HourVal
    $Sheet.Cells.Item($Row, 4).Value = $LocationVal
    $Sheet.Cells.Item($Row, 5).Value = $ValueVal
    $Sheet.Cells.Item($Row, 6).Value = $StatusVal
}

# --- ALERTS sheet: add row 20 ---
$wsAlerts = $wb.Worksheets.Item("ALERTS")
Add-LogRow $wsAlerts 20 "2026-02-01" "11:30:37" "11:00" "Living Room" "CRITICAL" "FALL_DETECTED"

# --- mmWave sheet: add rows 63-66 ---
$wsMmWave = $wb.Worksheets.Item("mmWave")

Add-LogRow $wsMmWave 63 "2026-02-01" "11:29:44" "11:00" "Living Room" "PRESENCE_DETECTED" "Active"
Add-LogRow $wsMmWave 64 "2026-02-01" "11:29:52" "11:00" "Living Room" "PRESENCE_DETECTED" "Active"
Add-LogRow $wsMmWave 65 "2026-02-01" "11:30:02" "11:00" "Living Room" "PRESENCE_DETECTED" "Active"
Add-LogRow $wsMmWave 66 "2026-02-01" "11:30:13" "11:00" "Living Room" "PRESENCE_DETECTED" "Active"
